$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

# Copy formatting (date style) from the cell above so the new row's A cell
# reuses the existing date-format style instead of minting a new one.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4163)

$ws.Cells.Item($row, 1).Value = 42619.893819444442
$ws.Cells.Item($row, 2).Value = 20
$ws.Cells.Item($row, 3).Value = 60
$ws.Cells.Item($row, 4).Value = 34
$ws.Cells.Item($row, 5).Value = 60
$ws.Cells.Item($row, 6).Value = 35
$ws.Cells.Item($row, 7).Value = 17791
$ws.Cells.Item($row, 8).Value = 21140
$ws.Cells.Item($row, 9).Value = 2284
$ws.Cells.Item($row, 10).Value = 356
$ws.Cells.Item($row, 11).Value = 200
$ws.Cells.Item($row, 12).Value = 18
$ws.Cells.Item($row, 13).Value = 10
$ws.Cells.Item($row, 14).Value = "Bag"
